$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), styled like the other header
# cells (bold, bordered, centered) by copying H1's format.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New data columns I (I0) and J (IF) for rows 2-13.
$dataI = @(1, 2, 3, 2, 1, 3, 6, 5, 5, 4, 1, 1)
$dataJ = @(6, 6, 5, 8, 5, 5, 8, 7, 6, 5, 3, 2)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
